$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Admin" row (row 2: Admin / admin@example.com / admin123) was removed from
# the Employees sheet, so every row below it shifts up by one.
$ws.Rows(2).Delete()

# Deleting the row does not automatically repair the worksheet's hyperlinks in
# this runtime, so the hyperlink objects (and the underlying mailto: targets)
# need to be rebuilt to point at the employees that now occupy rows 2-4.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:emp1@example.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:emp2@example.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:emp3@example.com")

# Re-applying the Hyperlink style keeps the cell formatting identical to how it
# was before (Hyperlinks.Add otherwise leaves the cells with a freshly created
# style record instead of reusing the existing one).
$ws.Range("B2:B4").Style = "Hyperlink"
